$wb = $excel.ActiveWorkbook

# --- Update the source metrics on the "Metrics" sheet (B2:B13). ---
# Every other changed cell in the workbook (the "today" sheet's B11:B22,
# E11:E22, F11:F22, plus the TODAY()-1 cell A1) is a formula that derives
# from these values (directly or transitively), so simply updating the
# raw numbers here and letting the engine recalc reproduces the rest of
# the diff automatically.
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value2  = 309299.93000000005
$wsMetrics.Range("B3").Value2  = 264768.26999999996
$wsMetrics.Range("B4").Value2  = 94493.069999999992
$wsMetrics.Range("B5").Value2  = 12659
$wsMetrics.Range("B6").Value2  = 5512007.040000001
$wsMetrics.Range("B7").Value2  = 4665121.2300000004
$wsMetrics.Range("B8").Value2  = 1626449.9500000002
$wsMetrics.Range("B9").Value2  = 215366
$wsMetrics.Range("B10").Value2 = 33977388.029999994
$wsMetrics.Range("B11").Value2 = 31940396.389999997
$wsMetrics.Range("B12").Value2 = 11908171.989999995
$wsMetrics.Range("B13").Value2 = 1312996

# --- Update the saved cursor/selection on each sheet. ---
# The "Metrics" sheet's stored selection moves to E21, and the "today"
# sheet's moves to H11. Activate "Metrics" first (so its selection is
# recorded), then activate "today" last so it ends up the workbook's
# selected/active tab again, matching the original file.
$wsMetrics.Activate()
$wsMetrics.Range("E21").Select()

$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate()
$wsToday.Range("H11").Select()
